# Updated cryptos list - apply new Price / Volume(1h) values
# (values scraped fresh; stored as text to match source inlineStr cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target columns keep their text formatting so values such as
# "42.079.26" or "8.21" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.079.26"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.283.35"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "318.41"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "100.85"
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "38.87"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").Value = "0.0900"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "8.21"
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "0.950"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "15.10"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "2.632.59"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "2.291.46"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "42.218.81"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "12.72"
$ws.Range("E21").Value = "  +28.44%  "
$ws.Range("D22").Value = "72.52"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "266.24"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "10.76"
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").Value = "22.39"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").Value = "37.12"
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("D31").Value = "165.44"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "6.03"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "0.0867"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("E35").Value = "  -12.17%  "
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").Value = "2.76"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "68.05"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("D45").Value = "91.14"
$ws.Range("E45").Value = "  -8.42%  "
$ws.Range("D46").Value = "114.57"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "11.86"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").Value = "78.69"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "8.92"
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("D50").Value = "1.605.54"
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("E51").Value = "  -2.50%  "
